$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level view tweaks ---
$wb.Windows.Item(1).TabRatio = 0.322

# --- Copy row 3 formatting down to rows 4-6 so the new data rows pick up the
#     same fonts/borders/alignment/wrap without inventing new style records ---
$ws.Range("A3:J3").Copy()
$ws.Range("A4:J4").PasteSpecial(-4122)
$ws.Range("A5:J5").PasteSpecial(-4122)
$ws.Range("A6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Clear the two stray formatted-but-empty cells from the old layout ---
$ws.Range("J5").ClearContents()
$ws.Range("H6").ClearContents()

# --- Row heights for the four data rows (3 teams worth of new content) ---
$ws.Rows.Item(3).RowHeight = 191.25
$ws.Rows.Item(4).RowHeight = 199.5
$ws.Rows.Item(5).RowHeight = 199.5
$ws.Rows.Item(6).RowHeight = 191.25

# --- Fill in the four team rows. Order chosen so that new shared strings land
#     in the same append order as the authoritative edit. ---
$ws.Range("C3").Value2 = "Ivan Dario Ruiz Bernal"
$ws.Range("B4").Value2 = "Equipo 2"
$ws.Range("B5").Value2 = "Equipo 3"
$ws.Range("B6").Value2 = "Equipo 4"
$ws.Range("D3").Value2 = "Se requiere verificar el funcionamiento del servicio tipo Soap suma`nProbar un servicio se debe a varios factores como el de garantizar que el servicio esté funcionando correctamente y devolviendo los datos esperados. Esto implica asegurarse de que la petición sea correcta y que la respuesta del servidor sea la adecuada.`n"
$ws.Range("E3").Value2 = "El alcance de las pruebas se limitan a la funcionalidad especifica del servicio realizando pruebas funcionales automatizadas, comprobando que las peticiones realizadas funcionen correctamente.`nEl alcance de las pruebas cubre la verificación de que se manejen adecuadamente situaciones de éxito y cuando es posible errores."
$ws.Range("D4").Value2 = "Se requiere verificar el funcionamiento del servicio tipo Soap cambio de divisas`nProbar un servicio se debe a varios factores como el de garantizar que el servicio esté funcionando correctamente y devolviendo los datos esperados. Esto implica asegurarse de que la petición sea correcta y que la respuesta del servidor sea la adecuada.`n"
$ws.Range("D5").Value2 = "Se requiere verificar el funcionamiento del servicio tipo Rest actualización de comentarios`nProbar un servicio se debe a varios factores como el de garantizar que el servicio esté funcionando correctamente y devolviendo los datos esperados. Esto implica asegurarse de que la petición sea correcta y que la respuesta del servidor sea la adecuada.`n"
$ws.Range("D6").Value2 = "Se requiere verificar el funcionamiento del servicio tipo Rest eliminar un post`nProbar un servicio se debe a varios factores como el de garantizar que el servicio esté funcionando correctamente y devolviendo los datos esperados. Esto implica asegurarse de que la petición sea correcta y que la respuesta del servidor sea la adecuada.`n"

# --- Remaining cells on the new rows reuse already-existing shared strings ---
$ws.Range("A4").Value2 = $ws.Range("A3").Value2
$ws.Range("A5").Value2 = $ws.Range("A3").Value2
$ws.Range("A6").Value2 = $ws.Range("A3").Value2
$ws.Range("C4").Value2 = $ws.Range("C3").Value2
$ws.Range("C5").Value2 = $ws.Range("C3").Value2
$ws.Range("C6").Value2 = $ws.Range("C3").Value2
$ws.Range("E4").Value2 = $ws.Range("E3").Value2
$ws.Range("E5").Value2 = $ws.Range("E3").Value2
$ws.Range("E6").Value2 = $ws.Range("E3").Value2
for ($r = 4; $r -le 6; $r++) {
    for ($c = 6; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item(3, $c).Value2
    }
}

# --- View state: the workbook was left scrolled to row 6 with E6 selected ---
$ws.Range("E6").Select()

